# Writing_Task_Evaluation_TG.xlsx — "Lots of minor changes trying to
# enforce class in RA-coded datasheets."
#
# 1. Row 46 (subject 229) had placeholder "-" text entries in C/D/E/G —
#    clear those cells entirely (blank, not a string placeholder).
# 2. Subject "278" appeared twice (rows 95 & 96) with no way to tell the
#    two raters apart. Relabel them "278 A" (row 95) and "278 B" (row 96)
#    so the Subject column disambiguates the duplicate.
# 3. The rating columns (B:G) were headed with bare letters A-F; rename
#    them rating1..rating6.
# 4. Turn the data range into a proper filterable table: AutoFilter over
#    A1:G125 (adds the _FilterDatabase defined name automatically).
# 5. Reset the selection to the header row of the renamed rating columns
#    and scroll down toward the bottom of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- clear the "-" placeholder cells in row 46 --------------------------
$ws.Range("C46").ClearContents()
$ws.Range("D46").ClearContents()
$ws.Range("E46").ClearContents()
$ws.Range("G46").ClearContents()

# --- disambiguate the duplicated "278" subject rows ----------------------
$ws.Range("A96").Value = "278 B"
$ws.Range("A95").Value = "278 A"

# --- rename the rating headers -------------------------------------------
$ws.Range("B1").Value = "rating1"
$ws.Range("C1").Value = "rating2"
$ws.Range("D1").Value = "rating3"
$ws.Range("E1").Value = "rating4"
$ws.Range("F1").Value = "rating5"
$ws.Range("G1").Value = "rating6"

# --- turn the range into an auto-filtered table ---------------------------
[void]$ws.Range("A1:G125").AutoFilter()
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$125")
$fd.Visible = $false

# --- selection / scroll state ---------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 50
[void]$ws.Range("B1:G1").Select()
